$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Developer name (was a "Student Name" placeholder)
$ws.Range("C3").Value = "Parneet kaur"

# Fill Preconditions (column E) for all 6 test cases first
$ws.Range("E7").Value = "No account exists. Create a new account with valid data."
$ws.Range("E8").Value = "Account created with invalid minimum balance type."
$ws.Range("E9").Value = "Account with balance=500.00, minimum_balance=100.00"
$ws.Range("E10").Value = "Account with balance=50.00, minimum_balance=50.00"
$ws.Range("E11").Value = "Account with balance=40.00, minimum_balance=50.00"
$ws.Range("E12").Value = "Account created with valid data."

# Then fill Method Inputs (column F) for all 6 test cases
$ws.Range("F7").Value = "account_number=9483914, client_number=22, balance=1000.00, date_created=`"2025-10-27`", minimum_balance=50.00"
$ws.Range("F8").Value = "minimum_balance=`"invalid`""
$ws.Range("F9").Value = "Call get_service_charges()"
$ws.Range("F10").Value = "Call get_service_charges()"
$ws.Range("F11").Value = "Call get_service_charges()"
$ws.Range("F12").Value = "Call str(account)"

# Then fill Expected Result (column G) for all 6 test cases
$ws.Range("G7").Value = "All attributes correctly set. Balance = 1000.00, Minimum Balance = 50.00"
$ws.Range("G8").Value = "minimum_balance defaults to 50.00"
$ws.Range("G9").Value = "Expected = 0.50"
$ws.Range("G10").Value = "Expected = 0.50"
$ws.Range("G11").Value = "Expected = 1.00"
$ws.Range("G12").Value = "Returns string: `"Account Number: 9483914 Balance: `$1,000.00\nMinimum Balance: `$50.00 Account Type: Savings`""

# Update selection to match the authored state (active cell G12)
$ws.Range("G12").Select()
